$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update first name: Robby -> Joe
$ws.Range("B2").Value = "Joe"

# Update last name: Jonas -> Jonas2
$ws.Range("C2").Value = "Jonas2"

# Update interests list (remove "Dance & Movement, " and swap order of the last two entries)
$ws.Range("G2").Value = "Education, Technology, Environment, Sports & Recreation, Coding & Software Development, Music & Performance, Health & Wellness, Animal Welfare"
